$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (instead of auto-converted numbers) for price cells
# whose new value looks numeric, matching the original inlineStr text cells.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values (prices, volume %, and the two coin swaps).
$ws.Range("D2").Value = '29.429.69'
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").Value = '1.837.90'
$ws.Range("E3").Value = '  -0.86%  '
$ws.Range("D4").Value = '0.9991'
$ws.Range("E4").Value = '  -1.04%  '
$ws.Range("D5").Value = '243.34'
$ws.Range("E5").Value = '  -0.75%  '
$ws.Range("D6").Value = '0.6260'
$ws.Range("E6").Value = '  +0.91%  '
$ws.Range("D7").Value = '1.0000'
$ws.Range("E7").Value = '  -1.50%  '
$ws.Range("D8").Value = '0.07403'
$ws.Range("E8").Value = '  -1.07%  '
$ws.Range("D9").Value = '0.2950'
$ws.Range("E9").Value = '  -0.50%  '
$ws.Range("D10").Value = '23.30'
$ws.Range("E10").Value = '  +0.78%  '
$ws.Range("E11").Value = '  -1.62%  '
$ws.Range("D12").Value = '1.838.51'
$ws.Range("E12").Value = '  -0.33%  '
$ws.Range("D13").Value = '5.009'
$ws.Range("E13").Value = '  -0.46%  '
$ws.Range("D14").Value = '0.6748'
$ws.Range("E14").Value = '  -0.33%  '
$ws.Range("D15").Value = '83.15'
$ws.Range("E15").Value = '  -0.37%  '
$ws.Range("D16").Value = '0.000009368'
$ws.Range("E16").Value = '  +3.10%  '
$ws.Range("D17").Value = '5.887'
$ws.Range("E17").Value = '  -0.54%  '
$ws.Range("D18").Value = '29.402.67'
$ws.Range("E18").Value = '  +0.23%  '
$ws.Range("D19").Value = '2.085.61'
$ws.Range("E19").Value = '  +0.23%  '
$ws.Range("D20").Value = '237.68'
$ws.Range("E20").Value = '  -0.53%  '
$ws.Range("D21").Value = '12.51'
$ws.Range("E21").Value = '  -1.72%  '
$ws.Range("E22").Value = '  -1.65%  '
$ws.Range("D23").Value = '7.350'
$ws.Range("E23").Value = '  +1.99%  '
$ws.Range("E24").Value = '  -1.38%  '
$ws.Range("D25").Value = '158.62'
$ws.Range("E25").Value = '  -1.34%  '
$ws.Range("D26").Value = '0.1417'
$ws.Range("E26").Value = '  -1.64%  '
$ws.Range("D27").Value = '8.471'
$ws.Range("E27").Value = '  -1.03%  '
$ws.Range("D28").Value = '17.74'
$ws.Range("E28").Value = '  -1.34%  '
$ws.Range("D29").Value = '0.06079'
$ws.Range("E29").Value = '  +7.79%  '
$ws.Range("D30").Value = '1.498'
$ws.Range("E30").Value = '  -0.86%  '
$ws.Range("D31").Value = '1.241'
$ws.Range("E31").Value = '  +1.26%  '
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").Value = '4.089'
$ws.Range("E32").Value = '  -1.03%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '4.106'
$ws.Range("E33").Value = '  -1.67%  '
$ws.Range("D34").Value = '1.842'
$ws.Range("E34").Value = '  -0.99%  '
$ws.Range("E35").Value = '  -0.63%  '
$ws.Range("D36").Value = '0.7252'
$ws.Range("E36").Value = '  -3.26%  '
$ws.Range("D37").Value = '2.609'
$ws.Range("E37").Value = '  -2.34%  '
$ws.Range("D38").Value = '2.885'
$ws.Range("E38").Value = '  +1.48%  '
$ws.Range("D39").Value = '1.219.94'
$ws.Range("E39").Value = '  -0.01%  '
$ws.Range("D40").Value = '0.01762'
$ws.Range("D41").Value = '6.296'
$ws.Range("E41").Value = '  -3.38%  '
$ws.Range("D42").Value = '0.9115'
$ws.Range("E42").Value = '  +0.40%  '
$ws.Range("E43").Value = '  -1.46%  '
$ws.Range("D44").Value = '1.998.95'
$ws.Range("E44").Value = '  +0.77%  '
$ws.Range("D45").Value = '101.79'
$ws.Range("E45").Value = '  -0.04%  '
$ws.Range("D46").Value = '65.42'
$ws.Range("E46").Value = '  -0.24%  '
$ws.Range("D47").Value = '0.5070'
$ws.Range("E47").Value = '  -1.82%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '9.251'
$ws.Range("E48").Value = '  +0.84%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.00000000118'
$ws.Range("E49").Value = '  -3.42%  '
$ws.Range("D50").Value = '0.4052'
$ws.Range("E50").Value = '  -0.47%  '
$ws.Range("D51").Value = '0.1137'
$ws.Range("E51").Value = '  +1.97%  '
